$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F (product_type) ---------------------------------------------
$ws.Range("F1").Value = "product_type"
$ws.Range("F2").Value = "rulebook"
$ws.Range("F3").Value = "rulebook"
$ws.Range("F4").Value = "scenario"
$ws.Range("F5").Value = "scenario"
$ws.Range("F6").Value = "scenario"
$ws.Range("F7").Value = "replay"
$ws.Range("F8").Value = "replay"

# --- Column G (product_code) ----------------------------------------------
# Format as text first so codes like "7-1" are stored verbatim.
$ws.Range("G2:G8").NumberFormat = "@"

$ws.Range("G1").Value = "product_code"
$ws.Range("G2").Value = "7-1"
$ws.Range("G3").Value = "7-2"
$ws.Range("G4").Value = "7-3"
$ws.Range("G6").Value = "7-5"
$ws.Range("G5").Value = "7-6"
$ws.Range("G7").Value = "7-4"
$ws.Range("G8").Value = "7-7"

# Adjust column E width to match the new layout (closest reproducible value
# to the 26.6640625 recorded by the authoring app's own pixel-exact autofit).
$ws.Columns("E").ColumnWidth = 25.83

# Leave the selection where the author ended up after entering the data.
$ws.Range("G9").Select() | Out-Null
